$wb = $excel.ActiveWorkbook

# --- Sheet "1": rename header G1 from "MVA-pris" to "MVA_cost", add J2 = H2 ---
$ws1 = $wb.Worksheets.Item("1")
$ws1.Range("G1").Value = "MVA_cost"
$ws1.Range("J2").Formula = "=H2"

# --- Sheet "2": same header rename, add J2 = H2 (formatted with 2 decimals) ---
$ws2 = $wb.Worksheets.Item("2")
$ws2.Range("G1").Value = "MVA_cost"
$ws2.Range("J2").Formula = "=H2"
$ws2.Range("J2").NumberFormat = "0.00"

# --- Costumers: add new "Kundenummer" column (G) ---
$ws5 = $wb.Worksheets.Item("Costumers")
$ws5.Range("G1").Value = "Kundenummer"
$ws5.Range("G2").Formula = "=A2+1"
$ws5.Columns.Item(7).AutoFit()

# --- Restore per-sheet selections, leaving Costumers as the active tab ---
$ws1.Activate()
$ws1.Range("J3").Select()

$ws2.Activate()
$ws2.Range("A2").Select()

$ws5.Activate()
$ws5.Range("G3").Select()
